$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Debug_Messages"
$ws.Range("B3").Value = "yes"
